# Thursday 03 June 2021 02:59:16 AM IST update
# Appends a new "Jun 2" snapshot (Daily/Weekly/Monthly/Closing RSI-screener
# columns) as 24 new trailing columns (CE:DB) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header labels for the 6 new 4-column blocks (CE1:DB1) ---
$ws.Cells.Item(1, 83).Value = 'Daily as on Jun 2, 16:00'
$ws.Cells.Item(1, 84).Value = 'Weekly as on Jun 2, 16:00'
$ws.Cells.Item(1, 85).Value = 'Monthly as on Jun 2, 16:00'
$ws.Cells.Item(1, 86).Value = 'Closing as on Jun 2, 16:00'
$ws.Cells.Item(1, 87).Value = 'Daily as on Jun 02, 16:00'
$ws.Cells.Item(1, 88).Value = 'Weekly as on Jun 02, 16:00'
$ws.Cells.Item(1, 89).Value = 'Monthly as on Jun 02, 16:00'
$ws.Cells.Item(1, 90).Value = 'Closing as on Jun 02, 16:00'
$ws.Cells.Item(1, 91).Value = 'Daily as on Jun 2, 15:51'
$ws.Cells.Item(1, 92).Value = 'Weekly as on Jun 2, 15:51'
$ws.Cells.Item(1, 93).Value = 'Monthly as on Jun 2, 15:51'
$ws.Cells.Item(1, 94).Value = 'Closing as on Jun 2, 15:51'
$ws.Cells.Item(1, 95).Value = 'Daily as on Jun 2, 15:56'
$ws.Cells.Item(1, 96).Value = 'Weekly as on Jun 2, 15:56'
$ws.Cells.Item(1, 97).Value = 'Monthly as on Jun 2, 15:56'
$ws.Cells.Item(1, 98).Value = 'Closing as on Jun 2, 15:56'
$ws.Cells.Item(1, 99).Value = 'Daily as on Jun 2, 15:58'
$ws.Cells.Item(1, 100).Value = 'Weekly as on Jun 2, 15:58'
$ws.Cells.Item(1, 101).Value = 'Monthly as on Jun 2, 15:58'
$ws.Cells.Item(1, 102).Value = 'Closing as on Jun 2, 15:58'
$ws.Cells.Item(1, 103).Value = 'Daily as on Jun 2, 15:49'
$ws.Cells.Item(1, 104).Value = 'Weekly as on Jun 2, 15:49'
$ws.Cells.Item(1, 105).Value = 'Monthly as on Jun 2, 15:49'
$ws.Cells.Item(1, 106).Value = 'Closing as on Jun 2, 15:49'

# Give the new header cells the same bold/border/centered look as the rest
# of row 1 (copy format from an existing header cell, e.g. B1).
$ws.Range("B1").Copy()
$ws.Range("CE1:DB1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows -----------------------------------------------------------
# Each stock row only got a new snapshot if it reappeared in that particular
# run, so the 4 new values land in whichever of the six new column blocks
# (CE:CH, CI:CL, CM:CP, CQ:CT, CU:CX, CY:DB) corresponds to that run.
# row 2 -> CE2:CH2
$ws.Cells.Item(2, 83).Value = 76.36
$ws.Cells.Item(2, 84).Value = 61.93
$ws.Cells.Item(2, 85).Value = 65.76
$ws.Cells.Item(2, 86).Value = 2201.25

# row 3 -> CE3:CH3
$ws.Cells.Item(3, 83).Value = 57.97
$ws.Cells.Item(3, 84).Value = 67.72
$ws.Cells.Item(3, 85).Value = 65.84
$ws.Cells.Item(3, 86).Value = 323

# row 4 -> CE4:CH4
$ws.Cells.Item(4, 83).Value = 64.26
$ws.Cells.Item(4, 84).Value = 63.36
$ws.Cells.Item(4, 85).Value = 68.4
$ws.Cells.Item(4, 86).Value = 5808.65

# row 5 -> CE5:CH5
$ws.Cells.Item(5, 83).Value = 67.42
$ws.Cells.Item(5, 84).Value = 72.68
$ws.Cells.Item(5, 85).Value = 70.27
$ws.Cells.Item(5, 86).Value = 11823.6

# row 6 -> CE6:CH6
$ws.Cells.Item(6, 83).Value = 57.88
$ws.Cells.Item(6, 84).Value = 59.44
$ws.Cells.Item(6, 85).Value = 63.38
$ws.Cells.Item(6, 86).Value = 1504

# row 7 -> CE7:CH7
$ws.Cells.Item(7, 83).Value = 52.98
$ws.Cells.Item(7, 84).Value = 52.33
$ws.Cells.Item(7, 85).Value = 61.72
$ws.Cells.Item(7, 86).Value = 675.9

# row 8 -> CE8:CH8
$ws.Cells.Item(8, 83).Value = 69.82
$ws.Cells.Item(8, 84).Value = 57.86
$ws.Cells.Item(8, 85).Value = 60.72
$ws.Cells.Item(8, 86).Value = 3014.6

# row 9 -> CE9:CH9
$ws.Cells.Item(9, 83).Value = 68.31
$ws.Cells.Item(9, 84).Value = 54.17
$ws.Cells.Item(9, 85).Value = 53.59
$ws.Cells.Item(9, 86).Value = 7184.7

# row 10 -> CE10:CH10
$ws.Cells.Item(10, 83).Value = 74.12
$ws.Cells.Item(10, 84).Value = 71.04
$ws.Cells.Item(10, 85).Value = 68.69
$ws.Cells.Item(10, 86).Value = 4295.05

# row 11 -> CE11:CH11
$ws.Cells.Item(11, 83).Value = 63.11
$ws.Cells.Item(11, 84).Value = 64.28
$ws.Cells.Item(11, 85).Value = 64.13
$ws.Cells.Item(11, 86).Value = 60.95

# row 12 -> CE12:CH12
$ws.Cells.Item(12, 83).Value = 57.74
$ws.Cells.Item(12, 84).Value = 62.29
$ws.Cells.Item(12, 85).Value = 53.82
$ws.Cells.Item(12, 86).Value = 342.3

# row 13 -> CE13:CH13
$ws.Cells.Item(13, 83).Value = 56.84
$ws.Cells.Item(13, 84).Value = 65.94
$ws.Cells.Item(13, 85).Value = 68.97
$ws.Cells.Item(13, 86).Value = 223.7

# row 14 -> CI14:CL14
$ws.Cells.Item(14, 87).Value = 44.46
$ws.Cells.Item(14, 88).Value = 64.22
$ws.Cells.Item(14, 89).Value = 61.99
$ws.Cells.Item(14, 90).Value = 2081.6

# row 15 -> CE15:CH15
$ws.Cells.Item(15, 83).Value = 45.52
$ws.Cells.Item(15, 84).Value = 49.13
$ws.Cells.Item(15, 85).Value = 58.03
$ws.Cells.Item(15, 86).Value = 529.5

# row 16 -> CE16:CH16
$ws.Cells.Item(16, 83).Value = 60.07
$ws.Cells.Item(16, 84).Value = 59.15
$ws.Cells.Item(16, 85).Value = 44.67
$ws.Cells.Item(16, 86).Value = 149.1

# row 17 -> CE17:CH17
$ws.Cells.Item(17, 83).Value = 72.55
$ws.Cells.Item(17, 84).Value = 76.58
$ws.Cells.Item(17, 85).Value = 83.81
$ws.Cells.Item(17, 86).Value = 543

# row 18 -> CE18:CH18
$ws.Cells.Item(18, 83).Value = 71.6
$ws.Cells.Item(18, 84).Value = 57.84
$ws.Cells.Item(18, 85).Value = 42.36
$ws.Cells.Item(18, 86).Value = 226.95

# row 19 -> CE19:CH19
$ws.Cells.Item(19, 83).Value = 54.1
$ws.Cells.Item(19, 84).Value = 44.61
$ws.Cells.Item(19, 85).Value = 39.92
$ws.Cells.Item(19, 86).Value = 8.8

# row 20 -> CQ20:CT20
$ws.Cells.Item(20, 95).Value = 63.95
$ws.Cells.Item(20, 96).Value = 59.58
$ws.Cells.Item(20, 97).Value = 50.84
$ws.Cells.Item(20, 98).Value = 1028

# row 21 -> CE21:CH21
$ws.Cells.Item(21, 83).Value = 69.93
$ws.Cells.Item(21, 84).Value = 66.73
$ws.Cells.Item(21, 85).Value = 51.38
$ws.Cells.Item(21, 86).Value = 110.85

# row 22 -> CU22:CX22
$ws.Cells.Item(22, 99).Value = 64.32
$ws.Cells.Item(22, 100).Value = 63.39
$ws.Cells.Item(22, 101).Value = 69.45
$ws.Cells.Item(22, 102).Value = 1917.95

# row 23 -> CE23:CH23
$ws.Cells.Item(23, 83).Value = 48.96
$ws.Cells.Item(23, 84).Value = 50.94
$ws.Cells.Item(23, 85).Value = 47.18
$ws.Cells.Item(23, 86).Value = 209

# row 24 -> CE24:CH24
$ws.Cells.Item(24, 83).Value = 63.38
$ws.Cells.Item(24, 84).Value = 62.63
$ws.Cells.Item(24, 85).Value = 60.76
$ws.Cells.Item(24, 86).Value = 1471.05

# row 25 -> CE25:CH25
$ws.Cells.Item(25, 83).Value = 50.32
$ws.Cells.Item(25, 84).Value = 45.75
$ws.Cells.Item(25, 85).Value = 47.37
$ws.Cells.Item(25, 86).Value = 162.85

# row 26 -> CE26:CH26
$ws.Cells.Item(26, 83).Value = 61.07
$ws.Cells.Item(26, 84).Value = 53.23
$ws.Cells.Item(26, 85).Value = 54.57
$ws.Cells.Item(26, 86).Value = 163.8

# row 27 -> CE27:CH27
$ws.Cells.Item(27, 83).Value = 76.68
$ws.Cells.Item(27, 84).Value = 65.44
$ws.Cells.Item(27, 85).Value = 67.33
$ws.Cells.Item(27, 86).Value = 1414.9

# row 28 -> CE28:CH28
$ws.Cells.Item(28, 83).Value = 63.66
$ws.Cells.Item(28, 84).Value = 57.48
$ws.Cells.Item(28, 85).Value = 56.57
$ws.Cells.Item(28, 86).Value = 121.6

# row 29 -> CE29:CH29
$ws.Cells.Item(29, 83).Value = 58.71
$ws.Cells.Item(29, 84).Value = 64.5
$ws.Cells.Item(29, 85).Value = 61.68
$ws.Cells.Item(29, 86).Value = 984.4

# row 30 -> CE30:CH30
$ws.Cells.Item(30, 83).Value = 80.82
$ws.Cells.Item(30, 84).Value = 74.68
$ws.Cells.Item(30, 85).Value = 67.58
$ws.Cells.Item(30, 86).Value = 437.25

# row 31 -> CE31:CH31
$ws.Cells.Item(31, 83).Value = 58.7
$ws.Cells.Item(31, 84).Value = 67.92
$ws.Cells.Item(31, 85).Value = 71.01
$ws.Cells.Item(31, 86).Value = 108.1

# row 32 -> CE32:CH32
$ws.Cells.Item(32, 83).Value = 51.33
$ws.Cells.Item(32, 84).Value = 56.36
$ws.Cells.Item(32, 85).Value = 67.94
$ws.Cells.Item(32, 86).Value = 3129.45

# row 33 -> CY33:DB33
$ws.Cells.Item(33, 103).Value = 49.11
$ws.Cells.Item(33, 104).Value = 58.83
$ws.Cells.Item(33, 105).Value = 60.14
$ws.Cells.Item(33, 106).Value = 711.85

# row 34 -> CM34:CP34
$ws.Cells.Item(34, 91).Value = 62.82
$ws.Cells.Item(34, 92).Value = 63.62
$ws.Cells.Item(34, 93).Value = 51.11
$ws.Cells.Item(34, 94).Value = 10.25

# row 35 -> CE35:CH35
$ws.Cells.Item(35, 83).Value = 65.4
$ws.Cells.Item(35, 84).Value = 68.26
$ws.Cells.Item(35, 85).Value = 71.91
$ws.Cells.Item(35, 86).Value = 2903.7

# row 36 -> CE36:CH36
$ws.Cells.Item(36, 83).Value = 64.26
$ws.Cells.Item(36, 84).Value = 69.5
$ws.Cells.Item(36, 85).Value = 81.2
$ws.Cells.Item(36, 86).Value = 811.5

# row 37 -> CE37:CH37
$ws.Cells.Item(37, 83).Value = 51.56
$ws.Cells.Item(37, 84).Value = 64.21
$ws.Cells.Item(37, 85).Value = 67.11
$ws.Cells.Item(37, 86).Value = 677.45

# row 38 -> CE38:CH38
$ws.Cells.Item(38, 83).Value = 57.88
$ws.Cells.Item(38, 84).Value = 77.89
$ws.Cells.Item(38, 85).Value = 79.44
$ws.Cells.Item(38, 86).Value = 1124.05

# row 39 -> CE39:CH39
$ws.Cells.Item(39, 83).Value = 49.64
$ws.Cells.Item(39, 84).Value = 69.35
$ws.Cells.Item(39, 85).Value = 84.32
$ws.Cells.Item(39, 86).Value = 1752.35

